# Apply updated statistics (2019-20 NCES data for free/reduced lunch and race
# variables) to the covariate_importance and strategy_importance sheets, and
# re-sort the rows by the new "positives" values, as described in the commit
# message.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: covariate_importance
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("covariate_importance")

# New row order (A=name, B=positives, C=ranks, D=pos_ranked), rows 2-3 unchanged.
$sheet1Data = @(
    @("percenttwoormoreraces", 97, 77, 77),
    @("percentasian", 82, 43, 43),
    @("percentwhite", 72, 37, 37),
    @("percentfreereducedlunch", 69, 35, 35),
    @("rplthemes", 62, 31, 31),
    @("percentblackorafricanamerican", 68, 23, 23),
    @("schoollevel", 60, 21, 21),
    @("derivedtotalenrolled", 47, 12, 12),
    @("percenthispaniclatino", 45, 8, 8),
    @("cntycaseschange", 56.99999999999999, 6, 6),
    @("percentamericanindianoralaskanative", 46, 6, 6),
    @("percentnativehawaiianorotherpacificislander", 16, 1, 1),
    @("locale", 28, 0, 0),
    @("percentnotspecified", 20, 0, 0)
)

$row = 4
foreach ($item in $sheet1Data) {
    $ws1.Cells.Item($row, 1).Value = $item[0]
    $ws1.Cells.Item($row, 2).Value = $item[1]
    $ws1.Cells.Item($row, 3).Value = $item[2]
    $ws1.Cells.Item($row, 4).Value = $item[3]
    $row++
}

# ---------------------------------------------------------------------------
# Sheet 2: strategy_importance
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("strategy_importance")

$sheet2Data = @(
    @("hvacsystems", 100, 100, 100),
    @("contacttracing", 97, 87, 94),
    @("cleaning", 91, 78, 90),
    @("screeningtestingforstudents", 89, 82, 89),
    @("masks", 68, 44, 66),
    @("physicaldistancing", 18, 4, 12),
    @("vaccination", 18, 4, 11),
    @("stayhome", 3, 0, 2),
    @("hepafilters", 5, 1, 1),
    @("quarantine", 2, 0, 1)
)

$row = 2
foreach ($item in $sheet2Data) {
    $ws2.Cells.Item($row, 1).Value = $item[0]
    $ws2.Cells.Item($row, 2).Value = $item[1]
    $ws2.Cells.Item($row, 3).Value = $item[2]
    $ws2.Cells.Item($row, 4).Value = $item[3]
    $row++
}
